$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.383.49"
$ws.Range("E2").Value = "  +1.31%  "

# Row 3
$ws.Range("D3").Value = "1.824.74"
$ws.Range("E3").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'313.91"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "'0.4469"
$ws.Range("E7").Value = "  +2.48%  "

# Row 8
$ws.Range("D8").Value = "'0.3756"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9
$ws.Range("D9").Value = "'0.07515"
$ws.Range("E9").Value = "  +2.70%  "

# Row 10
$ws.Range("E10").Value = "  +5.47%  "

# Row 11
$ws.Range("D11").Value = "'21.06"
$ws.Range("E11").Value = "  +1.66%  "

# Row 12
$ws.Range("D12").Value = "1.821.02"
$ws.Range("E12").Value = "  -0.32%  "

# Row 13
$ws.Range("D13").Value = "'6.764"
$ws.Range("E13").Value = "  +1.47%  "

# Row 14
$ws.Range("D14").Value = "'94.12"
$ws.Range("E14").Value = "  +5.07%  "

# Row 15
$ws.Range("D15").Value = "'5.418"
$ws.Range("E15").Value = "  +2.25%  "

# Row 16
$ws.Range("D16").Value = "'0.07112"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("E17").Value = "  +0.10%  "

# Row 18
$ws.Range("D18").Value = "'0.000008813"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("D20").Value = "'15.23"
$ws.Range("E20").Value = "  +2.18%  "

# Row 21
$ws.Range("D21").Value = "27.380.49"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22
$ws.Range("D22").Value = "'5.263"
$ws.Range("E22").Value = "  +2.30%  "

# Row 23
$ws.Range("D23").Value = "'10.93"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
$ws.Range("D24").Value = "2.057.24"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("D25").Value = "'1.974"
$ws.Range("E25").Value = "  -0.46%  "

# Row 26
$ws.Range("D26").Value = "'2.382"
$ws.Range("E26").Value = "  +7.44%  "

# Row 27
$ws.Range("D27").Value = "'151.44"
$ws.Range("E27").Value = "  +0.05%  "

# Row 28
$ws.Range("D28").Value = "'18.58"
$ws.Range("E28").Value = "  +1.58%  "

# Row 29
$ws.Range("D29").Value = "'5.363"
$ws.Range("E29").Value = "  +2.30%  "

# Row 30
$ws.Range("D30").Value = "'117.81"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31
$ws.Range("D31").Value = "'0.08835"
$ws.Range("E31").Value = "  +1.23%  "

# Row 32
$ws.Range("D32").Value = "'0.7867"
$ws.Range("E32").Value = "  +6.22%  "

# Row 33
$ws.Range("E33").Value = "  +1.87%  "

# Row 34
$ws.Range("D34").Value = "'4.523"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35
$ws.Range("D35").Value = "'2.906"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("E36").Value = "  +0.06%  "

# Row 37
$ws.Range("D37").Value = "'1.113"
$ws.Range("E37").Value = "  +1.70%  "

# Row 38
$ws.Range("D38").Value = "'0.01992"
$ws.Range("E38").Value = "  +2.52%  "

# Row 39
$ws.Range("D39").Value = "'0.05336"
$ws.Range("E39").Value = "  +1.94%  "

# Row 40
$ws.Range("D40").Value = "'7.396"
$ws.Range("E40").Value = "  +2.29%  "

# Row 41
$ws.Range("D41").Value = "'0.5316"
$ws.Range("E41").Value = "  +3.55%  "

# Row 42
$ws.Range("D42").Value = "'0.1731"
$ws.Range("E42").Value = "  +1.48%  "

# Row 43
$ws.Range("D43").Value = "'2.862"
$ws.Range("E43").Value = "  -0.12%  "

# Row 44
$ws.Range("D44").Value = "'2.287"
$ws.Range("E44").Value = "  +17.60%  "

# Row 45
$ws.Range("D45").Value = "'8.760"
$ws.Range("E45").Value = "  +1.97%  "

# Row 46
$ws.Range("D46").Value = "'0.5118"
$ws.Range("E46").Value = "  +7.44%  "

# Row 47
$ws.Range("D47").Value = "'10.65"
$ws.Range("E47").Value = "  +0.72%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.706"
$ws.Range("E48").Value = "  +2.61%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'106.06"
$ws.Range("E49").Value = "  +0.16%  "

# Row 50
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = "  +0.10%  "

# Row 51
$ws.Range("E51").Value = "  +0.59%  "
